$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-40 down to 35-41
$ws.Rows.Item(34).Insert()

# Fill in the new row 34 data
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44722
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = 100112035
$ws.Cells.Item(34, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 95
$ws.Cells.Item(34, 11).Value = 15000
$ws.Cells.Item(34, 12).Value = 15500
$ws.Cells.Item(34, 13).Value = 15263
$ws.Cells.Item(34, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 1018
$ws.Cells.Item(34, 17).Value = 15
$ws.Cells.Item(34, 18).Value = "Hortaliza"
